$wb = $excel.ActiveWorkbook

# Update the cell value from "qatitans40" to "qatitans65" on Sheet1 (A2)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "qatitans65"

# Make Sheet1 the active/selected sheet (instead of Sheet2)
$ws1.Activate()
